# Correción psp's Clase Configuración
#
# The underlying web-query report was refreshed/re-imported, which made
# Excel rename the sheet/connection/query-table from "excel" to "excel(1)"
# (and the matching defined name from "excel" to "excel_1"), the "generated
# at" timestamp string was refreshed, and several metric values in the
# report table changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (excel -> excel(1)) ------------------------------
$ws.Name = "excel(1)"

# --- Rename + re-point the workbook-scoped defined name -----------------
# (RefersTo auto-updates to the new sheet name once the sheet above was
# renamed, matching 'excel(1)'!$A$1:$D$29)
$name = $wb.Names.Item(1)
$name.Name = "excel_1"

# --- Refresh timestamp string in A27 -------------------------------------
$ws.Cells.Item(27, 1).Value2 = "Reporte generado a las 01:29 PM el 5/12/2018"

# --- Updated report metrics ----------------------------------------------
$ws.Range("C8").Value2  = 7.64
$ws.Range("D8").Value2  = 63.1
$ws.Range("D9").Value2  = 0.31458333333333333
$ws.Range("D15").Value2 = 2.99
$ws.Range("C17").Value2 = 105
$ws.Range("D17").Value2 = 628
